$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '66.911.77'
Set-TextValue 'E2' '  -2.38%  '
Set-TextValue 'D3' '2.662.08'
Set-TextValue 'E3' '  -1.23%  '
Set-TextValue 'E4' '  +0.11%  '
Set-TextValue 'D5' '594.07'
Set-TextValue 'E5' '  -0.80%  '
Set-TextValue 'D6' '163.92'
Set-TextValue 'E6' '  +2.77%  '
Set-TextValue 'E7' '  +0.02%  '
Set-TextValue 'D8' '0.544'
Set-TextValue 'E8' '  +0.22%  '
Set-TextValue 'D9' '2.662.60'
Set-TextValue 'E9' '  -1.14%  '
Set-TextValue 'E10' '  +0.73%  '
Set-TextValue 'E11' '  +0.67%  '
Set-TextValue 'E12' '  -0.64%  '
Set-TextValue 'E13' '  -1.94%  '
Set-TextValue 'D14' '27.65'
Set-TextValue 'E14' '  -2.37%  '
Set-TextValue 'D15' '3.150.92'
Set-TextValue 'E15' '  -1.82%  '
Set-TextValue 'D16' '0.0000181'
Set-TextValue 'E16' '  -2.87%  '
Set-TextValue 'D17' '67.056.39'
Set-TextValue 'E17' '  -1.99%  '
Set-TextValue 'D18' '2.660.13'
Set-TextValue 'E18' '  -1.13%  '
Set-TextValue 'E19' '  -2.53%  '
Set-TextValue 'D20' '359.84'
Set-TextValue 'E20' '  -1.77%  '
Set-TextValue 'D21' '7.47'
Set-TextValue 'E21' '  -1.27%  '
Set-TextValue 'D22' '4.36'
Set-TextValue 'E22' '  -3.61%  '
Set-TextValue 'D23' '4.78'
Set-TextValue 'E23' '  -3.14%  '
Set-TextValue 'E24' '  -5.36%  '
Set-TextValue 'E25' '  +0.21%  '
Set-TextValue 'D26' '71.10'
Set-TextValue 'E26' '  -5.18%  '
Set-TextValue 'D27' '10.02'
Set-TextValue 'E27' '  -1.50%  '
Set-TextValue 'D28' '2.800.68'
Set-TextValue 'E28' '  -1.46%  '
Set-TextValue 'B29' 'PEPE'
Set-TextValue 'C29' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue 'D29' '0.0000102'
Set-TextValue 'E29' '  -2.41%  '
Set-TextValue 'B30' 'Binance-PegBSC-USD'
Set-TextValue 'C30' 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue 'D30' '0.998'
Set-TextValue 'E30' '  -0.25%  '
Set-TextValue 'D31' '552.90'
Set-TextValue 'E31' '  -4.61%  '
Set-TextValue 'D32' '7.95'
Set-TextValue 'E32' '  -3.74%  '
Set-TextValue 'E33' '  -3.70%  '
Set-TextValue 'E34' '  -0.93%  '
Set-TextValue 'E35' '  -2.03%  '
Set-TextValue 'E36' '  +0.01%  '
Set-TextValue 'D37' '1.56'
Set-TextValue 'E37' '  -5.33%  '
Set-TextValue 'D38' '19.37'
Set-TextValue 'E38' '  -3.48%  '
Set-TextValue 'D39' '155.77'
Set-TextValue 'E39' '  -3.23%  '
Set-TextValue 'D40' '0.371'
Set-TextValue 'E40' '  -2.25%  '
Set-TextValue 'E41' '  -3.08%  '
Set-TextValue 'D42' '1.81'
Set-TextValue 'E42' '  -4.92%  '
Set-TextValue 'E43' '  +0.12%  '
Set-TextValue 'B44' 'USDe'
Set-TextValue 'C44' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D44' '1.00'
Set-TextValue 'E44' '  +0.02%  '
Set-TextValue 'B45' 'dogwifhat'
Set-TextValue 'C45' 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue 'D45' '2.52'
Set-TextValue 'E45' '  -4.93%  '
Set-TextValue 'D46' '40.21'
Set-TextValue 'E46' '  -0.62%  '
Set-TextValue 'D47' '0.0₆0297'
Set-TextValue 'E47' '  -6.04%  '
Set-TextValue 'D48' '0.583'
Set-TextValue 'E48' '  -2.99%  '
Set-TextValue 'D49' '152.24'
Set-TextValue 'E49' '  -3.98%  '
Set-TextValue 'E50' '  -2.55%  '
Set-TextValue 'E51' '  -2.83%  '
